$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to text format so numeric-looking values
# (e.g. "1.00", "464.50") are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "67.456.37"
$ws.Range("E2").Value = "  +1.08%  "

$ws.Range("D3").Value = "3.864.82"
$ws.Range("E3").Value = "  +0.24%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "464.50"
$ws.Range("E5").Value = "  +9.44%  "

$ws.Range("D6").Value = "148.63"
$ws.Range("E6").Value = "  +13.11%  "

$ws.Range("D7").Value = "0.634"
$ws.Range("E7").Value = "  +3.52%  "

$ws.Range("D8").Value = "0.998"
$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").Value = "0.753"
$ws.Range("E9").Value = "  +3.39%  "

$ws.Range("E10").Value = "  -2.76%  "

$ws.Range("D11").Value = "0.0000313"
$ws.Range("E11").Value = "  -8.79%  "

$ws.Range("D12").Value = "43.81"
$ws.Range("E12").Value = "  +6.96%  "

$ws.Range("D13").Value = "10.44"
$ws.Range("E13").Value = "  +1.56%  "

$ws.Range("D14").Value = "4.488.38"
$ws.Range("E14").Value = "  +0.27%  "

$ws.Range("D15").Value = "14.79"
$ws.Range("E15").Value = "  -7.03%  "

$ws.Range("D16").Value = "3.860.20"
$ws.Range("E16").Value = "  -0.24%  "

$ws.Range("E17").Value = "  -0.19%  "

$ws.Range("D18").Value = "20.12"
$ws.Range("E18").Value = "  +0.59%  "

$ws.Range("E19").Value = "  +7.65%  "

$ws.Range("D20").Value = "67.484.32"
$ws.Range("E20").Value = "  +0.69%  "

$ws.Range("D21").Value = "433.75"
$ws.Range("E21").Value = "  +4.93%  "

$ws.Range("E22").Value = "  -1.05%  "

$ws.Range("E23").Value = "  +8.07%  "

$ws.Range("D24").Value = "88.77"
$ws.Range("E24").Value = "  +5.20%  "

$ws.Range("E25").Value = "  +9.53%  "

$ws.Range("D26").Value = "10.39"
$ws.Range("E26").Value = "  +13.77%  "

$ws.Range("D27").Value = "37.63"
$ws.Range("E27").Value = "  -0.19%  "

$ws.Range("E28").Value = "  +1.83%  "

$ws.Range("E29").Value = "  +4.41%  "

$ws.Range("D30").Value = "743.71"
$ws.Range("E30").Value = "  +2.69%  "

$ws.Range("E31").Value = "  +10.55%  "

$ws.Range("E32").Value = "  +4.40%  "

$ws.Range("E33").Value = "  -0.80%  "

$ws.Range("D34").Value = "43.21"
$ws.Range("E34").Value = "  +10.29%  "

$ws.Range("D35").Value = "0.163"
$ws.Range("E35").Value = "  +6.59%  "

$ws.Range("D36").Value = "57.40"
$ws.Range("E36").Value = "  +3.45%  "

$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.16%  "

$ws.Range("D38").Value = "5.57"
$ws.Range("E38").Value = "  +2.34%  "

$ws.Range("D39").Value = "0.0480"
$ws.Range("E39").Value = "  +3.79%  "

$ws.Range("D40").Value = "0.352"
$ws.Range("E40").Value = "  +12.40%  "

$ws.Range("D41").Value = "2.64"
$ws.Range("E41").Value = "  +16.28%  "

$ws.Range("E42").Value = "  +0.35%  "

$ws.Range("E43").Value = "  +5.08%  "

$ws.Range("D44").Value = "0.0₃0678"
$ws.Range("E44").Value = "  -10.08%  "

$ws.Range("E45").Value = "  -0.18%  "

$ws.Range("E46").Value = "  +2.51%  "

$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "2.77"
$ws.Range("E47").Value = "  +7.80%  "

$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "3.25"
$ws.Range("E48").Value = "  +4.01%  "

$ws.Range("D49").Value = "2.14"
$ws.Range("E49").Value = "  +3.71%  "

$ws.Range("E50").Value = "  +3.28%  "

$ws.Range("D51").Value = "144.39"
$ws.Range("E51").Value = "  +1.59%  "

# Restore default style on column D so no stray formatting remains
$ws.Range("D2:D51").Style = "Normal"